$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final dataset for rows 2-31 (A:timestamp, B:label, C:ax, D:ay, E:az, F:gx, G:gy, H:gz)
# 5 new rows prepended (t=0..400), original series shifted to t=500..2400,
# and 5 new rows appended (t=2500..2900).
$rows = @(
  @(0, "falling", -2.278594136238097, 10.26830673217773, -1.992950439453123, 0.299476683139801, 0.1832595765590667, 0.4633412957191467),
  @(100, "falling", -1.48284924030304, 9.543427348136902, -0.6066013872623449, -0.0821613743901252, -0.3101668357849121, 0.4265366494655609),
  @(200, "falling", -2.347809791564941, 9.136160850524902, -1.348505258560181, -0.1818851232528686, 0.0901026204228401, 0.0786489024758338),
  @(300, "falling", -2.276926577091217, 9.172239780426027, -1.260899052023888, 0.0897971913218498, 0.1569923609495163, 0.2553416788578033),
  @(400, "falling", -3.074454665184022, 9.00393795967102, -1.032875627279282, 0.1241583600640297, 0.194713294506073, 0.1459967941045761),
  @(500, "falling", -2.568546533584593, 8.544089555740355, -1.202380612492562, 0.0265726372599601, -0.0716239511966705, -0.0607810914516449),
  @(600, "falling", -2.869771957397464, 8.188082218170166, -1.89351409673691, 0.2611449062824249, 0.041233405470848, 0.2226603776216507),
  @(700, "falling", -3.721651673316956, 7.623609662055968, -2.820488899946214, 0.2386955916881561, -0.3608686327934265, -0.0772744566202163),
  @(800, "falling", -4.097623348236085, 7.048315763473512, -3.382553160190582, 0.5322163701057434, -0.0166460778564214, 0.039248090237379),
  @(900, "falling", -4.406465649604798, 7.42472231388092, -2.858199380338192, 0.7162395119667053, -0.1989893466234207, -0.0884227454662323),
  @(1000, "falling", -1.668206214904804, 6.401832580566415, -0.706982076168021, 1.458746194839478, 0.5387831330299377, -0.3553708493709564),
  @(1100, "falling", -10.42951822280886, 10.53292512893678, 17.01628107577567, 1.215621829032898, 3.648850917816162, -3.939317226409912),
  @(1200, "falling", -13.75426650047303, 12.83328127861023, 8.825744509696944, 0.697149932384491, -3.561955213546753, 2.918103218078613),
  @(1300, "falling", -6.927797257900221, 11.75843620300293, 3.359905034303661, -1.915978789329529, 0.1592831164598465, -0.5018258094787598),
  @(1400, "falling", 0.7311251163482728, 8.272733926773066, 2.376347184181214, 1.59848165512085, 1.512196898460388, 0.6305656433105469),
  @(1500, "falling", 2.263805747032156, 7.363986670970927, 2.197247743606565, -0.159435823559761, -0.8156577944755554, 0.9199630618095398),
  @(1600, "falling", 0.7175407409668102, 9.555598974227896, 1.757197499275209, 0.0126754539087414, 0.1342376321554184, -1.328173756599426),
  @(1700, "falling", 3.773898720741262, 8.314922451972965, 2.525771677494061, -0.0963639914989471, -0.2640464901924133, -0.8072584271430969),
  @(1800, "falling", 1.089288711547856, 9.25660705566406, 5.577280521392808, 1.544267296791077, 0.0639881342649459, 0.7490735054016113),
  @(1900, "falling", 2.104704558849335, 8.927368879318237, 2.343123555183409, -1.995391368865967, 0.8813258409500122, -0.5760459303855896),
  @(2000, "falling", 2.28660249710083, 7.885151982307432, 2.352557301521302, -1.253037333488464, -0.7693848013877869, -0.4392121136188507),
  @(2100, "falling", 2.81174808740616, 8.390844881534578, 1.742006152868269, -0.1798998117446899, 0.3457497358322143, 0.0044287731871008),
  @(2200, "falling", 3.457506418228151, 8.694531917572021, 1.465195775032044, -0.0542142912745475, 0.2333505302667617, 0.1989893466234207),
  @(2300, "falling", 3.469298005104063, 8.534674406051636, 1.470959067344665, -0.0705549344420433, -0.119118720293045, -0.1950187236070633),
  @(2400, "falling", 2.979763031005859, 8.45093321800232, 1.210365951061249, 0.0215329993516206, 0.1395827084779739, -0.0713185146450996),
  @(2500, "falling", 3.198545038700105, 8.029186010360718, 1.33356249332428, 0.0612392425537109, 0.2464841306209564, 0.1278235465288162),
  @(2600, "falling", 3.5916051864624, 8.19729232788086, 1.436704874038696, 0.0074830991216003, -0.0374154970049858, -0.0755945742130279),
  @(2700, "falling", 3.165686726570128, 8.03831171989441, 1.441773623228073, 0.0713185146450996, 0.0108428578823804, 0.0226020142436027),
  @(2800, "falling", 3.061142683029175, 7.987302541732788, 1.411497831344604, 0.0126754539087414, 0.0655152946710586, 0.0435241498053073),
  @(2900, "falling", 3.059285700321198, 8.002557039260864, 1.58986583352089, 0.0337503030896186, -0.0120645882561802, -0.0161879286170005)
)

$data = New-Object 'object[,]' $rows.Count,8
for ($i = 0; $i -lt $rows.Count; $i++) {
  for ($j = 0; $j -lt 8; $j++) {
    $data[$i,$j] = $rows[$i][$j]
  }
}
$ws.Range("A2:H31").Value = $data